$d = $word.ActiveDocument

# p0_merge1
$d.Content.Find.Execute("Because we do not know when the user input will arrive, we should check often. We could check the serial receive buffer contiguously, but this would consume more power than checking it periodically", $false, $false, $false, $false, $false, $true, 1, $false, "Because we do not know when the user input will arrive, we should check often. We could check the serial receive buffer contiguously, but this would consume more power than checking it periodically", 2) | Out-Null

# p0_merge2
$d.Content.Find.Execute(". Checking periodically at a rate of one check per quarter second will be no different to the eyes of the user. We implemented this by creating a FreeRTOS task with a FreeRTOS delay of 250ms by dividing 250ms by the ", $false, $false, $false, $false, $false, $true, 1, $false, ". Checking periodically at a rate of one check per quarter second will be no different to the eyes of the user. We implemented this by creating a FreeRTOS task with a FreeRTOS delay of 250ms by dividing 250ms by the ", 2) | Out-Null

# p1_serialevent
$d.Content.Find.Execute("‘SerialEvent()", $false, $false, $false, $false, $false, $true, 1, $false, "‘SerialEvent()", 2) | Out-Null

# p1_loop
$d.Content.Find.Execute("‘loop()", $false, $false, $false, $false, $false, $true, 1, $false, "‘loop()", 2) | Out-Null

# p3_listen
$d.Content.Find.Execute("Due to not knowing when the first beacon will arrive, we have the same situation as with the user input. We could listen ", $false, $false, $false, $false, $false, $true, 1, $false, "Due to not knowing when the first beacon will arrive, we have the same situation as with the user input. We could listen ", 2) | Out-Null

# p3_consume
$d.Content.Find.Execute("consume more power. Therefore, we only process the beacons at a rate of one per quarter second. We did this with the help of FreeRTOS delay of 250ms as explained before. To be clear we do not disable the LoRa receiving during the FreeRTOS delay, because then the LoRa receive buffer would obviously not contain any received beacons.", $false, $false, $false, $false, $false, $true, 1, $false, "consume more power. Therefore, we only process the beacons at a rate of one per quarter second. We did this with the help of FreeRTOS delay of 250ms as explained before. To be clear we do not disable the LoRa receiving during the FreeRTOS delay, because then the LoRa receive buffer would obviously not contain any received beacons.", 2) | Out-Null

# p4_listeningforbeacons
$d.Content.Find.Execute("‘listeningForBeacons()", $false, $false, $false, $false, $false, $true, 1, $false, "‘listeningForBeacons()", 2) | Out-Null

# p6_vtasksuspendall
$d.Content.Find.Execute("‘vTaskSuspendAll()", $false, $false, $false, $false, $false, $true, 1, $false, "‘vTaskSuspendAll()", 2) | Out-Null

# p6_xtaskresumeall
$d.Content.Find.Execute("‘xTaskResumeAll()", $false, $false, $false, $false, $false, $true, 1, $false, "‘xTaskResumeAll()", 2) | Out-Null

# p8_input
$d.Content.Find.Execute(" We also experimented with setting the pins as input, but this resulted in an increase of current by 0.1mA.", $false, $false, $false, $false, $false, $true, 1, $false, " We also experimented with setting the pins as input, but this resulted in an increase of current by 0.1mA.", 2) | Out-Null

# p13_vapplicationidlehook
$d.Content.Find.Execute("‘vApplicationIdleHook()", $false, $false, $false, $false, $false, $true, 1, $false, "‘vApplicationIdleHook()", 2) | Out-Null

# p13_adcnoise
$d.Content.Find.Execute("‘ADC Noise Reduction’ mode. The ‘ADC Noise Reduction’ mode compared to the ‘Idle’ mode it consumes 0.7 mA less current. We also experimented with the ‘Power-save’ mode, the ‘standby’ mode and the ‘Extended standby’ mode but we could not get a reliable acknowledge back to the gateway in those modes. To sleep we first set the mode using ‘set", $false, $false, $false, $false, $false, $true, 1, $false, "‘ADC Noise Reduction’ mode. The ‘ADC Noise Reduction’ mode compared to the ‘Idle’ mode it consumes 0.7 mA less current. We also experimented with the ‘Power-save’ mode, the ‘standby’ mode and the ‘Extended standby’ mode but we could not get a reliable acknowledge back to the gateway in those modes. To sleep we first set the mode using ‘set", 2) | Out-Null

# p13_mode_to_the
$d.Content.Find.Execute("_mode’ to the ‘ADC Noise Reduction’ mode, then disable the interrupts to execute the following function without interrupts, set sleep to be enabled, then reenable the interrupts and finally going to sleep via the ‘sleep", $false, $false, $false, $false, $false, $true, 1, $false, "_mode’ to the ‘ADC Noise Reduction’ mode, then disable the interrupts to execute the following function without interrupts, set sleep to be enabled, then reenable the interrupts and finally going to sleep via the ‘sleep", 2) | Out-Null

# p13_powersave
$d.Content.Find.Execute("believe that the ‘Power-save’ mode would be the best to implement here since we only need to wake up after a sleep of around two to nine seconds.", $false, $false, $false, $false, $false, $true, 1, $false, "believe that the ‘Power-save’ mode would be the best to implement here since we only need to wake up after a sleep of around two to nine seconds.", 2) | Out-Null

# final_content
$d.Content.Find.Execute(" which disables the brown out detection. All the above settings results in a current usage of 11.3 mA, the same current usage as in the between the transmissions which leads us to conclude that we forgot to disable a core component of the PCB.`rTo find this last components we tried al lot of things, including a ‘power\_all\_diasble()’ function of the avr/pwer.h library which weirdly increased the power consumption.", $false, $false, $false, $false, $false, $true, 1, $false, " which disables the brown out detection. All the above settings results in a current usage of 11.3 mA, the same current usage as in the between the transmissions which prompted us to test the sleep modes. After various testing we found out that we in fact do not go into the sleep mode. We tried a lot of libraries including a popular one \url{https://www.arduino.cc/en/Reference/LowPowerDeepSleep} which is compatible with our Atmega32u4 but unfortunately it is not compatible with the FreeRTOS library.", 2) | Out-Null

